# Insert a new date column "30-dec" right before the existing "01-oct."
# column (FA) on the "Prix Spot" sheet, shifting all subsequent columns
# (FA..GE) one position to the right (FB..GF). The new column's header
# (row 1) gets the label "30-dec"; the new column's data rows (2-25) get
# the placeholder "-" used throughout the sheet for missing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at FA - this shifts FA:GE right to FB:GF and
# inherits formatting from the column being split (keeps header style).
$ws.Range("FA1").EntireColumn.Insert()

# New header label for the inserted column.
$ws.Range("FA1").Value = "30-dec"

# New column's data rows: fill with the sheet's standard "no data" marker.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 157).Value = "-"
}
